$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift columns D:K (4..11) to take the values previously in F:M (6..13) -- i.e. drop the
# two oldest quarters and shift everything left by two columns.
$rows = 8,9,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27
foreach ($r in $rows) {
    $old = @()
    for ($c = 4; $c -le 13; $c++) {
        $old += ,($ws.Cells.Item($r, $c).Value2)
    }
    for ($i = 0; $i -lt 8; $i++) {
        $ws.Cells.Item($r, 4+$i).Value2 = $old[$i+2]
    }
}

# Fill in the two new rightmost columns (L=12, M=13) with the newly-published quarters
$ws.Cells.Item(8, 12).Value2 = "12 ماهه منتهی به 1401/10"
$ws.Cells.Item(8, 13).Value2 = "3 ماهه منتهی به 1402/01"
$ws.Cells.Item(9, 12).Value2 = "1402-02-30 (3)"
$ws.Cells.Item(9, 13).Value2 = "1402-02-30"
$ws.Cells.Item(11, 12).Value2 = 19856
$ws.Cells.Item(11, 13).Value2 = 3377
$ws.Cells.Item(12, 12).Value2 = -12733
$ws.Cells.Item(12, 13).Value2 = -1960
$ws.Cells.Item(13, 12).Value2 = 7123
$ws.Cells.Item(13, 13).Value2 = 1417
$ws.Cells.Item(14, 12).Value2 = -2712
$ws.Cells.Item(14, 13).Value2 = -550
$ws.Cells.Item(15, 12).Value2 = "-"
$ws.Cells.Item(15, 13).Value2 = "-"
$ws.Cells.Item(16, 12).Value2 = -186
$ws.Cells.Item(16, 13).Value2 = -169
$ws.Cells.Item(17, 12).Value2 = 4225
$ws.Cells.Item(17, 13).Value2 = 698
$ws.Cells.Item(18, 12).Value2 = -126
$ws.Cells.Item(18, 13).Value2 = -28
$ws.Cells.Item(19, 12).Value2 = 492
$ws.Cells.Item(19, 13).Value2 = 46
$ws.Cells.Item(20, 12).Value2 = 4591
$ws.Cells.Item(20, 13).Value2 = 716
$ws.Cells.Item(21, 12).Value2 = -939
$ws.Cells.Item(21, 13).Value2 = -93
$ws.Cells.Item(22, 12).Value2 = 3653
$ws.Cells.Item(22, 13).Value2 = 623
$ws.Cells.Item(23, 12).Value2 = "-"
$ws.Cells.Item(23, 13).Value2 = "-"
$ws.Cells.Item(24, 12).Value2 = 3653
$ws.Cells.Item(24, 13).Value2 = 623
$ws.Cells.Item(25, 12).Value2 = 0
$ws.Cells.Item(25, 13).Value2 = 0
$ws.Cells.Item(26, 12).Value2 = 3782
$ws.Cells.Item(26, 13).Value2 = 2447
$ws.Cells.Item(27, 12).Value2 = 0
$ws.Cells.Item(27, 13).Value2 = 0

# A couple of cells in the newly shifted-in column I (9) needed a data correction
# (per the updated read_price algorithm / source filing corrections)
$ws.Cells.Item(19, 9).Value2 = "-"
$ws.Cells.Item(20, 9).Value2 = 740
$ws.Cells.Item(22, 9).Value2 = 696
$ws.Cells.Item(24, 9).Value2 = 696
